# Weekly update: a new price-report row for Berenjena (Vega Monumental
# Concepción) is inserted above the current row 8, pushing the existing
# rows 8-36 down to 9-37. The new row carries the same price/quality data
# as the (old) row 8, but is dated one reporting period later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts rows 8..36 down to 9..37)
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with this week's data
$ws.Cells.Item(8, 1).Value = 11
$ws.Cells.Item(8, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value = "Bíobío"
$ws.Cells.Item(8, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 100112001
$ws.Cells.Item(8, 7).Value = "Berenjena"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 9000
$ws.Cells.Item(8, 12).Value = 10000
$ws.Cells.Item(8, 13).Value = 9500
$ws.Cells.Item(8, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 158
$ws.Cells.Item(8, 17).Value = 60
$ws.Cells.Item(8, 18).Value = "Hortaliza"
